$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Try copying format from row 12 (Solved group) to row 9
$ws.Range("A12:H12").Copy()
$ws.Range("A9:H9").PasteSpecial(-4122)  # xlPasteFormats = -4122
